$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits within the "Profesor Asociado" block (row 2) ---
# Extend employment period to absorb the merged "Profesor Asistente" stint.
$ws.Range("B2").Value = "Ene. 2015 - Actualmente"

# Clarify committee membership bullets (now scoped to the Facultad de Psicología).
$ws.Range("E6").Value = "Miembro del Comité de Investigación de la Facultad de Psicología"
$ws.Range("E7").Value = "Miembro del Comité Asesor de Asuntos Éticos de la Facultad de Psicología"

# Merge the two doctoral co-supervision bullets into a single entry (row 15,
# before the row deletions below renumber things).
$ws.Range("E15").Value = "Co-supervisión de estudiantes de doctorado: \href{https://www.researchgate.net/profile/Milena-Vasquez-Amezquita}{Milena Vásquez-Amézquita} (PhD en Neurociencia, Universidad de Valencia, España - 2015-2018). Francisco Javier Flores  (Professional Doctorate in Counselling Psychology, University of East London, Reino Unido – 2016-2018)"

# --- Remove the now-redundant "Profesor Asistente" entry (rows 8-15) ---
# Its unique bullets (Organizador, Estancia posdoctoral, merged co-supervision)
# survive by deleting only the duplicated/obsolete rows, letting the rest
# shift up into the "Profesor Asociado" block above.
$ws.Rows(14).Delete()   # Co-supervisión: Milena (now folded into the merged bullet)
$ws.Rows(12).Delete()   # Supervisión de pregrado (duplicate)
$ws.Rows(10).Delete()   # Asesor metodológico (duplicate)
$ws.Rows(9).Delete()    # Líder del grupo CODEC (duplicate)
$ws.Rows(8).Delete()    # Profesor Asistente header row

# Match the saved selection state (whole row 11 selected).
$ws.Rows(11).Select()
